$d = $word.ActiveDocument

# 1. Replace the "Research (approx. 600 words)..." paragraph with the new Spring paragraph.
$d.Content.Find.Execute(
    "Research (approx. 600 words): Jakarta EE makes it possible to create enterprise applications. Carry out some research and critically discuss what other technologies and frameworks are available and how they compare to Jakarta EE",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Spring is a great technology to be used for developing the enterprise application, at the beginning could be much harder to be configured and maybe will take a while to use all the things incorporated in the framework but will be much easier than to create application because will give you the readability, the stability and the flexibility and the user experience will increase much quicker that way because of the flexible way of implement a JavaScript framework for user interface like react.js, angular or vue.js and validation of a form could be made asynchronous with to wait until a form is completed to give a hint is something is wrong or not.",
    2)

# 2. The three blank paragraphs that followed become the three new body paragraphs; the
#    final blank paragraph before the section break is left untouched.
$p = $d.Paragraphs.Item(4)
$p.Range.InsertBefore("A strong argument for using this technology for building enterprise applications consist in the fact showed in the statistics because in the development process using java as programming language most of the companies adopt to use this framework.")

$p = $d.Paragraphs.Item(5)
$p.Range.InsertBefore("Jakarta EE from my knowledge is not very flexible as the framework described above, for the user interface you cannot use whatever the JS library you want and will be much harder to configure and implement and it" + [char]0x2019 + "s not so used in production (Java EE).")

$p = $d.Paragraphs.Item(6)
$p.Range.InsertBefore("A future perspective is that the application will need to move to the cloud services and some of them will need to adapt to the new requirements or to be rewritten, this being a part of the process of the migration.")
